$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D values are stored as text (many look like ambiguous numbers e.g. "1.003")
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "26.095.45"
$ws.Range("E2").Value = "  -0.49%  "
$ws.Range("D3").Value = "1.648.43"
$ws.Range("E3").Value = "  -0.74%  "
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  -0.15%  "
$ws.Range("D5").Value = "218.26"
$ws.Range("E5").Value = "  -0.09%  "
$ws.Range("D6").Value = "0.5199"
$ws.Range("E6").Value = "  -0.76%  "
$ws.Range("D7").Value = "1.003"
$ws.Range("E7").Value = "  -0.19%  "
$ws.Range("D8").Value = "0.2617"
$ws.Range("E8").Value = "  -1.00%  "
$ws.Range("D9").Value = "0.06299"
$ws.Range("E9").Value = "  -0.36%  "
$ws.Range("D10").Value = "20.31"
$ws.Range("E10").Value = "  -1.89%  "
$ws.Range("D11").Value = "0.07659"
$ws.Range("E11").Value = "  -1.78%  "
$ws.Range("D12").Value = "4.585"
$ws.Range("E12").Value = "  +1.67%  "
$ws.Range("D13").Value = "1.679.45"
$ws.Range("E13").Value = "  +0.98%  "
$ws.Range("D14").Value = "1.873.24"
$ws.Range("E14").Value = "  -0.92%  "
$ws.Range("D15").Value = "0.5575"
$ws.Range("E15").Value = "  -0.97%  "
$ws.Range("D16").Value = "0.0₅8108"
$ws.Range("E16").Value = "  +0.54%  "
$ws.Range("D17").Value = "65.00"
$ws.Range("E17").Value = "  -0.41%  "
$ws.Range("D18").Value = "26.036.09"
$ws.Range("E18").Value = "  -0.76%  "
$ws.Range("E19").Value = "  +0.03%  "
$ws.Range("D20").Value = "4.594"
$ws.Range("E20").Value = "  -2.54%  "
$ws.Range("D21").Value = "193.67"
$ws.Range("E21").Value = "  -0.31%  "
$ws.Range("D22").Value = "10.43"
$ws.Range("E22").Value = "  +1.81%  "
$ws.Range("D23").Value = "5.921"
$ws.Range("E23").Value = "  -1.70%  "
$ws.Range("D24").Value = "1.004"
$ws.Range("E24").Value = "  -0.11%  "
$ws.Range("D25").Value = "144.59"
$ws.Range("E25").Value = "  -1.17%  "
$ws.Range("D26").Value = "0.1178"
$ws.Range("E26").Value = "  -2.97%  "
$ws.Range("D27").Value = "7.181"
$ws.Range("E27").Value = "  -0.85%  "
$ws.Range("D28").Value = "15.82"
$ws.Range("E28").Value = "  -1.83%  "
$ws.Range("D29").Value = "1.502"
$ws.Range("E29").Value = "  +1.04%  "
$ws.Range("D30").Value = "0.05404"
$ws.Range("E30").Value = "  -4.71%  "
$ws.Range("D31").Value = "1.268"
$ws.Range("E31").Value = "  -0.29%  "
$ws.Range("D32").Value = "3.439"
$ws.Range("E32").Value = "  -1.44%  "
$ws.Range("D33").Value = "3.323"
$ws.Range("E33").Value = "  -1.13%  "
$ws.Range("D34").Value = "1.554"
$ws.Range("E34").Value = "  -3.32%  "
$ws.Range("D35").Value = "2.414"
$ws.Range("E35").Value = "  +0.54%  "
$ws.Range("D36").Value = "2.780"
$ws.Range("E36").Value = "  -0.78%  "
$ws.Range("D37").Value = "0.9405"
$ws.Range("E37").Value = "  -0.35%  "
$ws.Range("D38").Value = "0.5572"
$ws.Range("E38").Value = "  -3.56%  "
$ws.Range("D39").Value = "0.01572"
$ws.Range("E39").Value = "  -1.67%  "
$ws.Range("B40").Value = "FraxShare"
$ws.Range("C40").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D40").Value = "5.745"
$ws.Range("E40").Value = "  -4.14%  "
$ws.Range("B41").Value = "PaxDollar"
$ws.Range("C41").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D41").Value = "1.003"
$ws.Range("E41").Value = "  -0.09%  "
$ws.Range("D42").Value = "1.027.24"
$ws.Range("E42").Value = "  -4.34%  "
$ws.Range("D43").Value = "0.8240"
$ws.Range("E43").Value = "  -3.01%  "
$ws.Range("D44").Value = "100.77"
$ws.Range("E44").Value = "  -2.11%  "
$ws.Range("D45").Value = "1.783.07"
$ws.Range("E45").Value = "  -1.00%  "
$ws.Range("E46").Value = "  +6.38%  "
$ws.Range("D47").Value = "57.19"
$ws.Range("E47").Value = "  -1.59%  "
$ws.Range("D48").Value = "0.9995"
$ws.Range("E48").Value = "  -0.44%  "
$ws.Range("D49").Value = "0.4314"
$ws.Range("E49").Value = "  -0.81%  "
$ws.Range("D50").Value = "7.876"
$ws.Range("E50").Value = "  -2.27%  "
$ws.Range("E51").Value = "  -4.08%  "
